$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$ws1 = $wb.Worksheets.Item("Metadata")

# Version: 2.0.0-sd-202406-matchbox-patch -> 2.0.1-sd-202510-matchbox-patch
$ws1.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"

# Date: 2024-06-19T17:47:42+02:00 -> 2025-10-29T22:15:57+01:00
$ws1.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new row for "Jurisdiction" after "Contact" (row 10), before "Description" (old row 11)
$ws1.Rows.Item(11).Insert()

# Match formatting of the surrounding data rows (copy format from the row that is now 12)
$ws1.Range("A12:B12").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)

$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""

# --- Elements sheet ---
$ws2 = $wb.Worksheets.Item("Elements")

# Row 5 = DataEnterer.typeId ; column AJ = Constraint(s)
$ws2.Range("AJ5").Value = "II-1:An II instance must have either a root or an nullFlavor. {root.exists() or nullFlavor.exists()}`n"
